# Natmi following Dr Hou advice
# Update LR-pair results (Ccl2 -> Ackr4) with the recalculated stats for the
# three clusters (ECs, FAPs, sCs), expanding the sending/target cluster
# combinations from 3 rows to all 6 sending x target pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @("ECs",  "Ccl2", "Ackr4", "FAPs", 3, 1, 56.47185866666666,  169.415576,        0.8070274173741353, 0.8070274173741354, 3, 1,                  1.890486333333333, 5.671459, 0.9442400689667343,  0.9442400689667344,  106.7592770272649,  960.833493245384,   0.762027624239399,   0.7620276242393992)
$data[3] = @("ECs",  "Ccl2", "Ackr4", "sCs",  3, 1, 56.47185866666666,  169.415576,        0.8070274173741353, 0.8070274173741354, 2, 0.6666666666666666, 0.1116383333333333, 0.334915, 0.05575993103326565, 0.05575993103326566, 6.304424181782222,  56.73981763604,     0.04499979313473628, 0.04499979313473629)
$data[4] = @("FAPs", "Ccl2", "Ackr4", "FAPs", 3, 1, 10.23495333333333,  30.70486,          0.1462655586439962, 0.1462655586439962, 3, 1,                  1.890486333333333, 5.671459, 0.9442400689667343,  0.9442400689667344,  19.34903939897111,  174.14135459074,    0.1381098011814649,  0.1381098011814649)
$data[5] = @("FAPs", "Ccl2", "Ackr4", "sCs",  3, 1, 10.23495333333333,  30.70486,          0.1462655586439962, 0.1462655586439962, 2, 0.6666666666666666, 0.1116383333333333, 0.334915, 0.05575993103326565, 0.05575993103326566, 1.142613131877778,  10.2835181869,      0.008155757462531299,0.008155757462531303)
$data[6] = @("sCs",  "Ccl2", "Ackr4", "FAPs", 3, 1, 3.268330666666667,  9.804992,          0.04670702398186845,0.04670702398186846,3, 1,                  1.890486333333333, 5.671459, 0.9442400689667343,  0.9442400689667344,  6.178734458147557,  55.60861012332801,  0.04410264354587038, 0.04410264354587039)
$data[7] = @("sCs",  "Ccl2", "Ackr4", "sCs",  3, 1, 3.268330666666667,  9.804992,          0.04670702398186845,0.04670702398186846,2, 0.6666666666666666, 0.1116383333333333, 0.334915, 0.05575993103326565, 0.05575993103326566, 0.3648709884088889, 3.28383889568,      0.00260438043599807, 0.00260438043599807)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($rowNum in 2..7) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $rowNum).Value = $values[$i]
    }
}
